$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.133160666666667
$ws.Range("H2").Value = 24.399482
$ws.Range("I2").Value = 0.3870696756706061
$ws.Range("J2").Value = 0.3870696756706061
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 413.6641385495631
$ws.Range("R2").Value = 3722.977246946068
$ws.Range("S2").Value = 0.1329319026224182
$ws.Range("T2").Value = 0.1329319026224183

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.133160666666667
$ws.Range("H3").Value = 24.399482
$ws.Range("I3").Value = 0.3870696756706061
$ws.Range("J3").Value = 0.3870696756706061
$ws.Range("M3").Value = 43.683024
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("Q3").Value = 355.281052597856
$ws.Range("R3").Value = 3197.529473380704
$ws.Range("S3").Value = 0.1141703664550796
$ws.Range("T3").Value = 0.1141703664550796

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.133160666666667
$ws.Range("H4").Value = 24.399482
$ws.Range("I4").Value = 0.3870696756706061
$ws.Range("J4").Value = 0.3870696756706061
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 298.0283322929771
$ws.Range("R4").Value = 2682.254990636794
$ws.Range("S4").Value = 0.09577207583428211
$ws.Range("T4").Value = 0.09577207583428213

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.133160666666667
$ws.Range("H5").Value = 24.399482
$ws.Range("I5").Value = 0.3870696756706061
$ws.Range("J5").Value = 0.3870696756706061
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 137.5292391487945
$ws.Range("R5").Value = 1237.76315233915
$ws.Range("S5").Value = 0.04419533075882612
$ws.Range("T5").Value = 0.04419533075882612

# Row 6
$ws.Range("I6").Value = 0.3559882250904906
$ws.Range("J6").Value = 0.3559882250904906
$ws.Range("M6").Value = 50.86142466666666
$ws.Range("N6").Value = 152.584274
$ws.Range("O6").Value = 0.3434314568613803
$ws.Range("P6").Value = 0.3434314568613804
$ws.Range("Q6").Value = 380.4471693906673
$ws.Range("R6").Value = 3424.024524516006
$ws.Range("S6").Value = 0.1222575547683242
$ws.Range("T6").Value = 0.1222575547683242

# Row 7
$ws.Range("I7").Value = 0.3559882250904906
$ws.Range("J7").Value = 0.3559882250904906
$ws.Range("M7").Value = 43.683024
$ws.Range("O7").Value = 0.294960761928139
$ws.Range("P7").Value = 0.294960761928139
$ws.Range("Q7").Value = 326.752208380752
$ws.Range("R7").Value = 2940.769875426768
$ws.Range("S7").Value = 0.105002558110137
$ws.Range("T7").Value = 0.105002558110137

# Row 8
$ws.Range("I8").Value = 0.3559882250904906
$ws.Range("J8").Value = 0.3559882250904906
$ws.Range("M8").Value = 36.64360566666667
$ws.Range("N8").Value = 109.930817
$ws.Range("O8").Value = 0.2474285170192034
$ws.Range("P8").Value = 0.2474285170192035
$ws.Range("Q8").Value = 274.0968453698803
$ws.Range("R8").Value = 2466.871608328923
$ws.Range("S8").Value = 0.08808163861043848
$ws.Range("T8").Value = 0.08808163861043847

# Row 9
$ws.Range("I9").Value = 0.3559882250904906
$ws.Range("J9").Value = 0.3559882250904906
$ws.Range("M9").Value = 16.90969166666667
$ws.Range("N9").Value = 50.729075
$ws.Range("O9").Value = 0.1141792641912772
$ws.Range("P9").Value = 0.1141792641912772
$ws.Range("Q9").Value = 126.4857280741583
$ws.Range("R9").Value = 1138.371552667425
$ws.Range("S9").Value = 0.04064647360159098
$ws.Range("T9").Value = 0.04064647360159098

# Row 10
$ws.Range("G10").Value = 5.398902333333333
$ws.Range("H10").Value = 16.196707
$ws.Range("I10").Value = 0.2569420992389033
$ws.Range("J10").Value = 0.2569420992389034
$ws.Range("M10").Value = 50.86142466666666
$ws.Range("N10").Value = 152.584274
$ws.Range("O10").Value = 0.3434314568613803
$ws.Range("P10").Value = 0.3434314568613804
$ws.Range("Q10").Value = 274.5958643095242
$ws.Range("R10").Value = 2471.362778785718
$ws.Range("S10").Value = 0.08824199947063793
$ws.Range("T10").Value = 0.08824199947063796

# Row 11
$ws.Range("G11").Value = 5.398902333333333
$ws.Range("H11").Value = 16.196707
$ws.Range("I11").Value = 0.2569420992389033
$ws.Range("J11").Value = 0.2569420992389034
$ws.Range("M11").Value = 43.683024
$ws.Range("O11").Value = 0.294960761928139
$ws.Range("P11").Value = 0.294960761928139
$ws.Range("Q11").Value = 235.840380200656
$ws.Range("R11").Value = 2122.563421805904
$ws.Range("S11").Value = 0.07578783736292243
$ws.Range("T11").Value = 0.07578783736292245

# Row 12
$ws.Range("G12").Value = 5.398902333333333
$ws.Range("H12").Value = 16.196707
$ws.Range("I12").Value = 0.2569420992389033
$ws.Range("J12").Value = 0.2569420992389034
$ws.Range("M12").Value = 36.64360566666667
$ws.Range("N12").Value = 109.930817
$ws.Range("O12").Value = 0.2474285170192034
$ws.Range("P12").Value = 0.2474285170192035
$ws.Range("Q12").Value = 197.8352481355132
$ws.Range("R12").Value = 1780.517233219619
$ws.Range("S12").Value = 0.06357480257448285
$ws.Range("T12").Value = 0.06357480257448286

# Row 13
$ws.Range("G13").Value = 5.398902333333333
$ws.Range("H13").Value = 16.196707
$ws.Range("I13").Value = 0.2569420992389033
$ws.Range("J13").Value = 0.2569420992389034
$ws.Range("M13").Value = 16.90969166666667
$ws.Range("N13").Value = 50.729075
$ws.Range("O13").Value = 0.1141792641912772
$ws.Range("P13").Value = 0.1141792641912772
$ws.Range("Q13").Value = 91.29377379511389
$ws.Range("R13").Value = 821.643964156025
$ws.Range("S13").Value = 0.02933745983086011
$ws.Range("T13").Value = 0.02933745983086012
